$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("AF3").Value = 67
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.38

# Row 4
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8

# Row 7
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7

# Row 8
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 3
$ws.Range("Q8").Value = 2.15
$ws.Range("R8").Value = 1.67

# Row 11
$ws.Range("I11").Value = 2.3
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("U11").Value = 2.1
$ws.Range("V11").Value = 1.67
$ws.Range("W11").Value = 8
$ws.Range("Z11").Value = 34
$ws.Range("AC11").Value = 6.5
$ws.Range("AH11").Value = 10
$ws.Range("AU11").Value = 9

# Row 12
$ws.Range("G12").Value = 1.9
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 4.75
$ws.Range("K12").Value = 1.95
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 6.5
$ws.Range("S12").Value = 1.57
$ws.Range("T12").Value = 2.25
$ws.Range("U12").Value = 2.2
$ws.Range("V12").Value = 1.62
$ws.Range("Z12").Value = 15
$ws.Range("AA12").Value = 19
$ws.Range("AC12").Value = 6.5
$ws.Range("AE12").Value = 19
$ws.Range("AG12").Value = 9.5
$ws.Range("AO12").Value = 11
$ws.Range("AP12").Value = 26
$ws.Range("AR12").Value = 67
$ws.Range("AS12").Value = 251
$ws.Range("AT12").Value = 2.25

# Row 13
$ws.Range("G13").Value = 2.6
$ws.Range("H13").Value = 2.8
$ws.Range("J13").Value = 3.6
$ws.Range("S13").Value = 1.73
$ws.Range("T13").Value = 2
$ws.Range("U13").Value = 2.5
$ws.Range("V13").Value = 1.5
$ws.Range("W13").Value = 5.5
$ws.Range("X13").Value = 11
$ws.Range("Y13").Value = 12
$ws.Range("AB13").Value = 51
$ws.Range("AC13").Value = 5
$ws.Range("AE13").Value = 23
$ws.Range("AF13").Value = 101
$ws.Range("AG13").Value = 6
$ws.Range("AO13").Value = 19
$ws.Range("AP13").Value = 41
$ws.Range("AQ13").Value = 67
$ws.Range("AR13").Value = 126
$ws.Range("AT13").Value = 2
$ws.Range("AU13").Value = 11
$ws.Range("AV13").Value = 101
$ws.Range("AZ13").Value = 81
$ws.Range("BA13").Value = 151
